$d = $word.ActiveDocument

# 1) Delete the large block of paragraphs that is dropped entirely in the
#    new review (old paragraphs 6 through 22 inclusive - everything between
#    the "מנגנון self-attention כגרף" paragraph and the final URL paragraph).
$startPara = $d.Paragraphs.Item(6)
$endPara = $d.Paragraphs.Item(22)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# 2) Paragraph 1: update the date and the paper title (keeps the <w:br/>
#    line break between the two runs).
$d.Content.Find.Execute("22.12.24", $true, $false, $false, $false, $false, `
    $true, 1, $false, "20.12.24", 2)
$d.Content.Find.Execute("Reasoning in Large Language Models: A Geometric Perspective", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "FAN: Fourier Analysis Networks", 2)

# 3) Paragraph 2: replace the whole abstract paragraph text.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "היום סוקרים קצרות מאמר המציע שכבה ארכיטקטונית חדשה לרשתות נוירונים. שכבה זו משלבת פונקציות מחזוריות כמו סינוס וקוסינוס. פונקציות מחזוריות אינן חיה חדשה בטריטוריה של הרשתות; כבר ראינו אותם במאמרי Neural radiance fields או NERF שהן משמשים לבניית מודלי 3D של אובייקטים וסצנות. למיטב זכרוני היה מאמר שבנה ייצוג של תמונה באמצעות רשת המערבת אקטיבציות מחזוריות. "

# 4) Paragraph 3: replace "רעיונות מרכזיים:" with the new paragraph text.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "אולם המאמר של היום מציע לבנות שכבה המכילה פונקציות מחזוריות אלא מציע לשלב אותן עם פונקציות אקטיבציות קלאסיות יותר כמו סיגמויד כאשר השילוב הוא לינארי. אז השכבה בגדול בנויה מצירוף לינארי של סינוסים וקוסינוסים עם מקדמים נלמדים יחד עם פונקציות אקטיבציות סטנדרטיות. השכבה הזו טובה למידול פונקציות מחזוריות כאשר ביצועיה על פונקציות לא מחזוריות אינן ברורות (המאמר טוען שיש שיפור גם שם), "

# 5) Paragraph 4: collapse the two runs (title + br + body) into one run.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "המאמר גם מציע להחליף ב-FAN את שכבות ה-FFN בטרנספורמרים וגם שכבות gating ב-LSTM (אותו סכום ממשוקל את סינוסים וקוסינוס יחד עם הסיגמואיד) ומדווח שיפור בביצועים בכמה משימות."

# 6) Paragraph 5: collapse the two runs into the short "רעיון מעניין…" line.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "רעיון מעניין…"

# 7) Final paragraph (now index 6 after the big deletion): update the link.
$d.Content.Find.Execute("https://arxiv.org/abs/2407.02678", $true, $false, $false, `
    $false, $false, $true, 1, $false, "https://arxiv.org/abs/2410.02675", 2)
